# DEV-5032: add support for overruling default permissions for classes/properties
# Adds a new "default_permissions_overrule" column (M) to the "classes" sheet,
# with explanatory comment and two example values, plus a doc hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classes")

# --- New header cell M1 -------------------------------------------------
$ws.Range("M1").Value = "default_permissions_overrule"

# Link the new header to the docs, like the other header cells (A1, B1, G1, L1)
$wb.Hyperlinks.Add($ws.Range("M1"), "https://docs.dasch.swiss/latest/DSP-TOOLS/file-formats/json-project/ontologies/", "default_permissions_overrule") | Out-Null

# Match the bold+underlined "header link" look used by the other linked headers.
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").Font.Underline = $true

# Explanation comment on the new header cell
$commentText = "optional`n`n" + `
  "If you have set this project to ""public"" in the json_header.xlsx, then you can still hide certain classes.`n`n" + `
  "Use ""private"" on any resource class (including image classes) to make it invisible for people outside of your project.`n`n" + `
  "Use ""limited view"" on an image class to blur the image for people outside of your project."
$ws.Range("M1").AddComment($commentText) | Out-Null

# --- Example values in the data rows ------------------------------------
$ws.Range("M5").Value = "limited_view"
$ws.Range("M6").Value = "private"

# --- Column sizing --------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = 26.5

# --- Selection, as left by the author after editing ----------------------
$ws.Range("P13").Select() | Out-Null
